# Updated solution for Tutorial 6
# Applies updated attendance dates (slashes -> dashes) and recalculated
# attendance counts for roll 2001ME72.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new date text (dashes instead of slashes).
$dates = @{
    3  = "28-07-2022"
    4  = "01-08-2022"
    5  = "04-08-2022"
    6  = "08-08-2022"
    7  = "11-08-2022"
    8  = "15-08-2022"
    9  = "18-08-2022"
    10 = "22-08-2022"
    11 = "25-08-2022"
    12 = "29-08-2022"
    13 = "01-09-2022"
    14 = "05-09-2022"
    15 = "08-09-2022"
    16 = "12-09-2022"
    17 = "15-09-2022"
    18 = "19-09-2022"
    19 = "22-09-2022"
    20 = "26-09-2022"
    21 = "29-09-2022"
}

foreach ($row in $dates.Keys) {
    $cell = $ws.Cells.Item($row, 1)
    # Some of the new day-month-year strings (day <= 12) look like valid
    # dates to Excel's smart-parsing, which would otherwise convert the
    # text into a date serial number. Force a text format for the write,
    # then clear formatting again so the cell stays visually identical to
    # before (no explicit style), while keeping the text content.
    $cell.NumberFormat = "@"
    $cell.Value = $dates[$row]
    $cell.ClearFormats()
}

# Map of row number -> new values for columns D (4), E (5), F (6), G (7), H (8).
# Only rows whose counts actually changed are listed here.
$counts = @{
    3  = @(1, 0, 0, 1, 1)
    4  = @(2, 1, 1, 0, 0)
    5  = @(1, 1, 0, 0, 0)
    6  = @(1, 1, 0, 0, 0)
    10 = @(1, 1, 0, 0, 0)
    11 = @(1, 0, 0, 1, 1)
    16 = @(1, 1, 0, 0, 0)
}

foreach ($row in $counts.Keys) {
    $vals = $counts[$row]
    $ws.Cells.Item($row, 4).Value = $vals[0]
    $ws.Cells.Item($row, 5).Value = $vals[1]
    $ws.Cells.Item($row, 6).Value = $vals[2]
    $ws.Cells.Item($row, 7).Value = $vals[3]
    $ws.Cells.Item($row, 8).Value = $vals[4]
}
